$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NCT(2.511289113398383, 1.3486432352281819, -0.5118117728483328, 2.6363948226455776)"
$ws.Range("C2").Value = "NIG(0.9444823312261025, 0.7519503641150205, 3.964639888727289, 5.793692820531108)"
$ws.Range("D2").Value = "NCT(2.992262237553456, 1.489604146397591, -0.8470358811827754, 2.3104717216487742)"
$ws.Range("E2").Value = "NIG(1.2361789815788045, 0.8477144170167379, 5.36855162485538, 5.7259158785691655)"
